{"js": "// Bold + color (\"2C3E50\") quantitative metrics (percentages, dollar\n// amounts, etc.) inside specific resume bullet paragraphs, matching the\n// target diff. Each target paragraph is located by its exact original\n// text, then the numeric/metric substrings inside it are located with a\n// paragraph-scoped search() and given bold + color run formatting -\n// Word automatically splits the run(s) around the match.\n\nconst HIGHLIGHT_COLOR = \"#2C3E50\";\n\n// Paragraph original text -> ordered list of metric substrings to bold+color.\nconst targets = [\n  {\n    text: \"\\u2022 Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%\",\n    metrics: [\"23%\", \"64%\"]\n  },\n  {\n    text: \"\\u2022 Utilized advanced sampling methods to decrease survey margin of error from \\u00B14.2% to \\u00B12.1%, increasing voter turnout prediction accuracy from 71% to 87%, and ensuring survey results more closely reflected true population attitudes\",\n    metrics: [\"\\u00B14.2%\", \"\\u00B12.1%\", \"71%\", \"87%\"]\n  },\n  {\n    text: \"\\u2022 Trigonometric algorithm for boundary estimation reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M and enabling smaller nonprofits to conduct analysis\",\n    metrics: [\"73.5%\", \"$4.7M\"]\n  },\n  {\n    text: \"\\u2022 Built real-time FEC analysis systems using Python, Pandas and PySpark to detect likely fraud, money laundering and financial crimes across billions of records daily, performing time series analysis on trillions of records in the political spending sub-economy valued over $2 trillion\",\n    metrics: [\"$2\"]\n  },\n  {\n    text: \"\\u2022 Predictive excellence: Utilized advanced sampling methods to decrease survey margin of error from \\u00B14.2% to \\u00B12.1%\",\n    metrics: [\"\\u00B14.2%\", \"\\u00B12.1%\"]\n  },\n  {\n    text: \"\\u2022 Increased voter turnout prediction accuracy from 71% to 87%\",\n    metrics: [\"71%\", \"87%\"]\n  },\n  {\n    text: \"\\u2022 Methodological advancement: Improved segmentation accuracy 34% and survey incidence 28%\",\n    metrics: [\"34%\", \"28%\"]\n  }\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  paragraphs.items[i].load(\"text\");\n}\nawait context.sync();\n\n// Map each target's exact text to the matching Paragraph object(s).\nfor (const target of targets) {\n  const match = paragraphs.items.find((p) => p.text === target.text);\n  if (!match) continue;\n\n  for (const metric of target.metrics) {\n    const found = match.search(metric, { matchCase: true });\n    found.load(\"items\");\n    await context.sync();\n\n    for (let i = 0; i < found.items.length; i++) {\n      found.items[i].font.bold = true;\n      found.items[i].font.color = HIGHLIGHT_COLOR;\n    }\n    await context.sync();\n  }\n}\n", "ps1": "# Bold + color (#2C3E50) quantitative metrics (percentages, dollar amounts,\n# etc.) inside specific resume bullet paragraphs, matching the target diff.\n# Each target paragraph is located by its exact original text; the metric\n# substrings inside it are then located with Find.Execute on a duplicated\n# sub-range and given bold + color run formatting -- Word automatically\n# splits the run(s) around the match, same as it would interactively.\n\n$HighlightColor = '2C3E50'\n\n$d = $word.ActiveDocument\n\n$targets = @(\n    @{ Text = '\u2022 Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%'; Metrics = @('23%', '64%') },\n    @{ Text = '\u2022 Utilized advanced sampling methods to decrease survey margin of error from \u00b14.2% to \u00b12.1%, increasing voter turnout prediction accuracy from 71% to 87%, and ensuring survey results more closely reflected true population attitudes'; Metrics = @('\u00b14.2%', '\u00b12.1%', '71%', '87%') },\n    @{ Text = '\u2022 Trigonometric algorithm for boundary estimation reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M and enabling smaller nonprofits to conduct analysis'; Metrics = @('73.5%', '$4.7M') },\n    @{ Text = '\u2022 Built real-time FEC analysis systems using Python, Pandas and PySpark to detect likely fraud, money laundering and financial crimes across billions of records daily, performing time series analysis on trillions of records in the political spending sub-economy valued over $2 trillion'; Metrics = @('$2') },\n    @{ Text = '\u2022 Predictive excellence: Utilized advanced sampling methods to decrease survey margin of error from \u00b14.2% to \u00b12.1%'; Metrics = @('\u00b14.2%', '\u00b12.1%') },\n    @{ Text = '\u2022 Increased voter turnout prediction accuracy from 71% to 87%'; Metrics = @('71%', '87%') },\n    @{ Text = '\u2022 Methodological advancement: Improved segmentation accuracy 34% and survey incidence 28%'; Metrics = @('34%', '28%') }\n)\n\n$paraCount = $d.Paragraphs.Count\n\nforeach ($target in $targets) {\n    for ($i = 1; $i -le $paraCount; $i++) {\n        $p = $d.Paragraphs($i)\n        $fullText = $p.Range.Text.TrimEnd(\"`r\")\n        if ($fullText -ne $target.Text) { continue }\n\n        $paraRange = $p.Range\n        foreach ($metric in $target.Metrics) {\n            $searchRange = $paraRange.Duplicate\n            $find = $searchRange.Find\n            $find.ClearFormatting()\n            $found = $find.Execute($metric)\n            if ($found) {\n                $searchRange.Font.Bold = 1\n                $searchRange.Font.Color = $HighlightColor\n            }\n        }\n        break\n    }\n}\n"}
